# Split the Title and Abstract paragraphs' single runs into one run per
# "word"/space token, mirroring the target OOXML diff. No formatting is
# changed -- only how the text is partitioned across <w:r> elements.

$d = $word.ActiveDocument

function Escape-Xml([string]$text) {
    $text = $text -replace '&', '&amp;'
    $text = $text -replace '<', '&lt;'
    $text = $text -replace '>', '&gt;'
    return $text
}

function Split-IntoTokens([string]$text) {
    # Alternating word / single-space tokens, e.g. "a b c" -> "a", " ", "b", " ", "c"
    $words = $text -split ' '
    $tokens = @()
    for ($i = 0; $i -lt $words.Count; $i++) {
        $tokens += $words[$i]
        if ($i -lt $words.Count - 1) {
            $tokens += " "
        }
    }
    return $tokens
}

function Set-ParagraphRuns($paragraph, [string]$pStyle, [string[]]$tokens) {
    $runsXml = ""
    foreach ($tok in $tokens) {
        $runsXml += '<w:r><w:t xml:space="preserve">' + (Escape-Xml $tok) + '</w:t></w:r>'
    }

    $pPrXml = ""
    if ($pStyle) {
        $pPrXml = '<w:pPr><w:pStyle w:val="' + $pStyle + '"/></w:pPr>'
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $pPrXml + $runsXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $paragraph.Range.InsertXML($xml)
}

function Find-ParagraphByText($doc, [string]$wantedText) {
    foreach ($para in $doc.Paragraphs) {
        $t = $para.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $wantedText) {
            return $para
        }
    }
    return $null
}

# --- Title paragraph: "Answers: Trigonometry (radians)" ---
$titleText = "Answers: Trigonometry (radians)"
$titlePara = Find-ParagraphByText $d $titleText
if ($titlePara -eq $null) { $titlePara = $d.Paragraphs(1) }
$titleTokens = @("Answers:", " ", "Trigonometry", " ", "(radians)")
Set-ParagraphRuns $titlePara "Title" $titleTokens

# --- Abstract paragraph: "Answers to the questions on trigonometry, using radians to measure angles." ---
$abstractText = "Answers to the questions on trigonometry, using radians to measure angles."
$abstractPara = Find-ParagraphByText $d $abstractText
if ($abstractPara -eq $null) { $abstractPara = $d.Paragraphs(3) }
$abstractTokens = Split-IntoTokens $abstractText
Set-ParagraphRuns $abstractPara "Abstract" $abstractTokens
